$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.520.31'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.917.14'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4824'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.69%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2893'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06722'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '111.88'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '18.93'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.924.64'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07556'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.321'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6711'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '299.29'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.538.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.00'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.580'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9998'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007578'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.167.67'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9994'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.447'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.481'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '165.14'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.36'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.70%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.110'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1065'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.438'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.150'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.066'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05018'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.98%  '
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('E35').Value = '  -1.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9998'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.738'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02024'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.32%  '
$ws.Range('E39').Value = '  -0.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '110.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.19%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.023'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4443'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8624'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +4.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.851'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9997'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.256'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.11%  '
$ws.Range('E48').Value = '  +1.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.286'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1238'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2520'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.75%  '
